# Add the three new character styles (GaNStyle, GaNParagraph, GaNLinks)
# and apply them to the runs that were newly styled in the commit.

$d = $word.ActiveDocument

# --- 1. Define the character styles -----------------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- 2. Apply GaNStyle to every "2022: Datumi kampanje ..." run --------

$datumiText = "2022: Datumi kampanje za opazovanje Ozvezdje škornjev: 14.-23. maj, 13.-22. junij, 12.-21. julij"
$r = $d.Content
$r.Find.ClearFormatting()
$guard = 0
while ($r.Find.Execute($datumiText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r.Style = "GaNStyle"
    $r.Collapse(0)
    $guard = $guard + 1
    if ($guard -gt 20) { break }
}

# --- 3. Apply GaNParagraph to the "Sodelujete v svetovni ..." run ------

$sodelujeteText = "Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega Ozvezdje škornjev na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom."
$r2 = $d.Content
$r2.Find.ClearFormatting()
if ($r2.Find.Execute($sodelujeteText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r2.Style = "GaNParagraph"
}

# --- 4. Apply GaNLinks to the "Jenik Hollan, CzechGlobe ..." run -------

$jenikText = "Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$r3 = $d.Content
$r3.Find.ClearFormatting()
if ($r3.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r3.Style = "GaNLinks"
}

Write-Output "Done"
